# Updated DPM integration testfixture with hierarchy node labels
$wb = $excel.ActiveWorkbook

# --- Sheet "CodeSchemes" ---
# Column A width target (OOXML chars): 28.6 -> 33.0
$wsCodeSchemes = $wb.Worksheets.Item("CodeSchemes")
$wsCodeSchemes.Range("A1").EntireColumn.ColumnWidth = 32.285714285714285
$wsCodeSchemes.Range("A2").Value = "85528a72-d49c-4255-abf6-83e10776f926"

# --- Sheet "Codes" ---
# Column A width target (OOXML chars): 35.2 -> 34.1
$wsCodes = $wb.Worksheets.Item("Codes")
$wsCodes.Range("A1").EntireColumn.ColumnWidth = 33.42857142857143
$wsCodes.Range("A2").Value = "56843db2-1331-4f96-bfd1-9dbb4aa417d0"
$wsCodes.Range("A3").Value = "f52af210-65e8-4395-8332-b8ab2588644b"

# --- Sheet "Extensions" ---
$wsExtensions = $wb.Worksheets.Item("Extensions")
$wsExtensions.Range("A2").Value = "823b3b57-9e89-43a6-9090-fbc7740e4f3e"

# --- Sheet "Members_dpmDimension" ---
# Column A width target (OOXML chars): 30.800000000000004 -> 35.2
$wsMembers = $wb.Worksheets.Item("Members_dpmDimension")
$wsMembers.Range("A1").EntireColumn.ColumnWidth = 34.42857142857143
$wsMembers.Range("A2").Value = "da6fb479-c860-4e73-942b-0867434d32ca"
$wsMembers.Range("A3").Value = "6b66a84d-979e-4ae4-94b9-2174add97d27"
